# Bookings workbook update:
#  - remove workbook protection flag
#  - add 5 new booking rows (as literal text, matching the source data entry)
#  - move selection / active cell to F9 and make the sheet tab active

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Workbook was unprotected (workbookProtection removed from workbook.xml).
$wb.Unprotect()

$bookings = @(
    @("Nikhil Sunny", "101",  "2025-07-21", "21:49", "1"),
    @("Nik",          "101",  "2025-07-21", "21:00", "1"),
    @("Nik",          "101",  "2025-07-21", "21:00", "1"),
    @("Chakku",       "N410", "2025-07-22", "09:00", "09:30"),
    @("Chakku",       "N410", "2025-07-22", "10:00", "10:30")
)

$row = 2
foreach ($booking in $bookings) {
    $rowRange = $ws.Range("A$row" + ":E$row")
    # Force plain text storage so values like "101" / "2025-07-21" / "09:30"
    # stay literal strings instead of being coerced to numbers/dates.
    $rowRange.NumberFormat = "@"

    $ws.Range("A$row").Value = $booking[0]
    $ws.Range("B$row").Value = $booking[1]
    $ws.Range("C$row").Value = $booking[2]
    $ws.Range("D$row").Value = $booking[3]
    $ws.Range("E$row").Value = $booking[4]

    $row++
}

# Active cell / selection moves to F9 and the sheet tab becomes selected.
$ws.Range("F9").Select()
